$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 863 (which holds the
# "Lluteño" / 2021-08-17 triplet). Everything currently at row 863
# and below shifts down by 3 rows, so the former rows 863-905 become
# rows 866-908 and the sheet's used range grows from A1:R905 to A1:R908.
$ws.Rows("863:865").Insert()

# Populate the 3 newly-inserted rows with a fresh Primera/Segunda/Tercera
# triplet for Choclo - Lluteño dated 2023-01-13 (serial 44939).

# Row 863 - Primera
$ws.Cells.Item(863, 1).Value = 1
$ws.Cells.Item(863, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(863, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(863, 4).Value = 44939
$ws.Cells.Item(863, 5).Value = 15
$ws.Cells.Item(863, 6).Value = 100112024
$ws.Cells.Item(863, 7).Value = "Choclo"
$ws.Cells.Item(863, 8).Value = "Lluteño"
$ws.Cells.Item(863, 9).Value = "Primera"
$ws.Cells.Item(863, 10).Value = 50
$ws.Cells.Item(863, 11).Value = 29000
$ws.Cells.Item(863, 12).Value = 30000
$ws.Cells.Item(863, 13).Value = 29500
$ws.Cells.Item(863, 14).Value = "`$/saco 50 unidades"
$ws.Cells.Item(863, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(863, 16).Value = 590
$ws.Cells.Item(863, 17).Value = 50
$ws.Cells.Item(863, 18).Value = "Hortaliza"

# Row 864 - Segunda
$ws.Cells.Item(864, 1).Value = 1
$ws.Cells.Item(864, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(864, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(864, 4).Value = 44939
$ws.Cells.Item(864, 5).Value = 15
$ws.Cells.Item(864, 6).Value = 100112024
$ws.Cells.Item(864, 7).Value = "Choclo"
$ws.Cells.Item(864, 8).Value = "Lluteño"
$ws.Cells.Item(864, 9).Value = "Segunda"
$ws.Cells.Item(864, 10).Value = 50
$ws.Cells.Item(864, 11).Value = 25000
$ws.Cells.Item(864, 12).Value = 26000
$ws.Cells.Item(864, 13).Value = 25500
$ws.Cells.Item(864, 14).Value = "`$/saco 75 unidades"
$ws.Cells.Item(864, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(864, 16).Value = 340
$ws.Cells.Item(864, 17).Value = 75
$ws.Cells.Item(864, 18).Value = "Hortaliza"

# Row 865 - Tercera
$ws.Cells.Item(865, 1).Value = 1
$ws.Cells.Item(865, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(865, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(865, 4).Value = 44939
$ws.Cells.Item(865, 5).Value = 15
$ws.Cells.Item(865, 6).Value = 100112024
$ws.Cells.Item(865, 7).Value = "Choclo"
$ws.Cells.Item(865, 8).Value = "Lluteño"
$ws.Cells.Item(865, 9).Value = "Tercera"
$ws.Cells.Item(865, 10).Value = 50
$ws.Cells.Item(865, 11).Value = 20000
$ws.Cells.Item(865, 12).Value = 21000
$ws.Cells.Item(865, 13).Value = 20500
$ws.Cells.Item(865, 14).Value = "`$/saco 100 unidades"
$ws.Cells.Item(865, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(865, 16).Value = 205
$ws.Cells.Item(865, 17).Value = 100
$ws.Cells.Item(865, 18).Value = "Hortaliza"
